$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Munka1")

# Rename model id in A24 from "BIOMD0000000991A" to "BIOMD0000000991"
$ws.Range("A24").Value = "BIOMD0000000991"

# Reflect the updated scroll/selection state left by the author when saving
$ws.Range("A25").Select()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
